$wb = $excel.ActiveWorkbook

# Rename sheet 2
$ws2 = $wb.Worksheets.Item("SearchModelDataByCondition")
$ws2.Name = "getConceptModelDataByCondition"

# Activate sheet 1 (testDataForMethod1) so it becomes the selected/active tab
$ws1 = $wb.Worksheets.Item("testDataForMethod1")
$ws1.Activate()
